# Master_Fuel_Sector_List.xlsx edit
# Adds a new sector row ("2L1_Oil-tanker-loading") to the "Sectors" sheet,
# inserted immediately above the existing "2L_Other-process-emissions" row,
# with activity=pop, units=1000, type=NC (same pattern as surrounding rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Insert a new row at row 43 (pushes 2L_Other-process-emissions and everything
# below it down by one row).
$ws.Rows.Item(43).Insert()

$ws.Cells.Item(43, 1).Value = "2L1_Oil-tanker-loading"
$ws.Cells.Item(43, 2).Value = "pop"
$ws.Cells.Item(43, 3).Value = 1000
$ws.Cells.Item(43, 4).Value = "NC"
